$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_strategy")

$ws.Range("A6").Value = "Strategy 321"
$ws.Range("A7").Value = "Strategy 4"
